# style pour la page admin et page changer les variables
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header sub-row (row 7): the tariff label that used to sit under "Avril"
# (P7, text "28") now belongs under "Mai" (Q7) and is corrected to "21";
# the placeholder "0" that used to sit in Q7 moves into P7 instead.
$ws.Range("T7").Copy()
$ws.Range("Q7").PasteSpecial(-4163)   # xlPasteValues: Q7 <- "28" (same text as T7)
$ws.Range("E7").Copy()
$ws.Range("P7").PasteSpecial(-4163)   # xlPasteValues: P7 <- "0"

# Correct "28" -> "21" everywhere it is shown (Q7 and T7). Build the text
# once in a scratch cell (forcing text via a leading apostrophe), then copy
# just its value into the real cells so their existing number formatting /
# style stays untouched.
$ws.Range("Z1").Value = "'21"
$ws.Range("Z1").Copy()
$ws.Range("Q7").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("T7").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("Z1").Clear()

# Monthly hours for Mai (column Q) updated for each stagiaire (rows 8-14).
$ws.Range("Q8").Value = 7
$ws.Range("Q9").Value = 7
$ws.Range("Q10").Value = 7
$ws.Range("Q11").Value = 7
$ws.Range("Q12").Value = 7
$ws.Range("Q13").Value = 7
$ws.Range("Q14").Value = 7

# Totals row for Mai (Q15) bumped accordingly.
$ws.Range("Q15").Value = 71

$ws.Application.CutCopyMode = 0
